# Auto-generated edit script applying the diff changes to before.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')
$ws.Range('F3').Value = 2138
$ws.Range('F5').Value = 1566
$ws.Range('F6').Value = 153
$ws.Range('F7').Value = 218
$ws.Range('F8').Value = 375
$ws.Range('F9').Value = 238
$ws.Range('F10').Value = 1156
$ws.Range('F11').Value = 683
$ws.Range('F12').Value = 457
$ws.Range('F13').Value = 750
$ws.Range('F14').Value = 78
$ws.Range('F15').Value = 223
$ws.Range('F16').Value = 183
$ws.Range('F17').Value = 246
$ws.Range('F18').Value = 154
$ws.Range('F19').Value = 293
$ws.Range('F20').Value = 1508
$ws.Range('F21').Value = 139
$ws.Range('F22').Value = 72
$ws.Range('F23').Value = 11
$ws.Range('F24').Value = 52
$ws.Range('F25').Value = 2217
$ws.Range('F26').Value = 115
$ws.Range('F27').Value = 778
$ws.Range('F28').Value = 53
$ws.Range('F29').Value = 75
$ws.Range('F30').Value = 57

$ws = $wb.Worksheets.Item('演出')
$ws.Range('F3').Value = 26
$ws.Range('F8').Value = 15
$ws.Range('F12').Value = 34
$ws.Range('F15').Value = 449
$ws.Range('F16').Value = 168
$ws.Range('F17').Value = 10
$ws.Range('F18').Value = 133

$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F2').Value = 1668
$ws.Range('F5').Value = 1771
$ws.Range('F6').Value = 1779
$ws.Range('F7').Value = 596
$ws.Range('F8').Value = 583
$ws.Range('F9').Value = 459

$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F2').Value = 1668
$ws.Range('B4').Value = '2023.10.25'
$ws.Range('C4').Value = '上海·方块大战（豫园店）'
$ws.Range('D4').Value = '丽水路88号2楼213 城隍庙第一购物中心'
$ws.Range('E4').Value = '2023.10.25 10:00-2024.10.20 21:00'
$ws.Range('F4').Value = 28
$ws.Range('G4').Value = '49'
$ws.Range('I4').Value = 'https://show.bilibili.com/platform/detail.html?id=79057&msource=Msearch_colligation'
$ws.Range('B5').Value = '2023.11.02'
$ws.Range('C5').Value = '上海·Hello Kitty Cosmos 50周年光影特展'
$ws.Range('D5').Value = '漕宝路3055号 宝龙美术馆'
$ws.Range('E5').Value = '2023.11.02 10:00-2024.01.28 18:00'
$ws.Range('F5').Value = 105
$ws.Range('G5').Value = '139'
$ws.Range('I5').Value = 'https://show.bilibili.com/platform/detail.html?id=77862&msource=Msearch_colligation'
$ws.Range('B6').Value = '2023.12.01'
$ws.Range('C6').Value = '上海·2023《蔚蓝档案》x  萌果酱谷子咖啡'
$ws.Range('D6').Value = '南京东路340号百联ZX 萌果酱谷子咖啡（百联）'
$ws.Range('E6').Value = '2023.12.01 00:00-2024.01.31 23:59'
$ws.Range('F6').Value = 1771
$ws.Range('I6').Value = 'https://show.bilibili.com/platform/detail.html?id=79005&msource=Msearch_colligation'
$ws.Range('B7').Value = '2023.12.06'
$ws.Range('C7').Value = '上海·「咒术回战  × animate cafe」'
$ws.Range('D7').Value = '西藏北路198号大悦城北座8楼N809-1 animate cafe上海店'
$ws.Range('E7').Value = '2023.12.06 00:00-2024.01.24 23:59'
$ws.Range('F7').Value = 1779
$ws.Range('I7').Value = 'https://show.bilibili.com/platform/detail.html?id=79292&msource=Msearch_colligation'
$ws.Range('B8').Value = '2023.12.09'
$ws.Range('C8').Value = '上海·非人哉官方授权主题店'
$ws.Range('D8').Value = '南京东路830号第一百货商业中心B馆5楼(海底捞旁边) 第一百货商业中心'
$ws.Range('E8').Value = '2023.12.09 00:00-2024.01.22 23:59'
$ws.Range('F8').Value = 596
$ws.Range('G8').Value = '30'
$ws.Range('I8').Value = 'https://show.bilibili.com/platform/detail.html?id=79240&msource=Msearch_colligation'
$ws.Range('B9').Value = '2023.12.10'
$ws.Range('C9').Value = '上海·多维跃迁-2023 红点设计概念大奖获奖作品展'
$ws.Range('D9').Value = '国展路1099号 上海世博展览馆'
$ws.Range('E9').Value = '2023.12.10 12:00-2024.02.16 17:00'
$ws.Range('F9').Value = 24
$ws.Range('G9').Value = '80'
$ws.Range('I9').Value = 'https://show.bilibili.com/platform/detail.html?id=78809&msource=Msearch_colligation'
$ws.Range('B10').Value = '2023.12.22'
$ws.Range('C10').Value = '上海·新海诚导演作品《铃芽之旅》展 丨 购票抽新海诚见面会门票丨 超限定复刻原画发售'
$ws.Range('D10').Value = '湖滨路168号 上海无限极荟购物中心'
$ws.Range('E10').Value = '2023.12.22 10:00-2024.02.16 22:00'
$ws.Range('F10').Value = 1566
$ws.Range('G10').Value = '65'
$ws.Range('I10').Value = 'https://show.bilibili.com/platform/detail.html?id=79166&msource=Msearch_colligation'
$ws.Range('F11').Value = 26
$ws.Range('F12').Value = 583
$ws.Range('F14').Value = 459
$ws.Range('F15').Value = 153
$ws.Range('F16').Value = 218
$ws.Range('F17').Value = 375
$ws.Range('F18').Value = 238
$ws.Range('F19').Value = 1156
$ws.Range('F20').Value = 683
$ws.Range('F21').Value = 457
$ws.Range('F23').Value = 15
$ws.Range('F24').Value = 750
$ws.Range('F25').Value = 78
$ws.Range('F26').Value = 223
$ws.Range('F29').Value = 34
$ws.Range('F30').Value = 183
$ws.Range('F31').Value = 246
$ws.Range('F32').Value = 154
$ws.Range('F33').Value = 293
$ws.Range('F35').Value = 1508
$ws.Range('F36').Value = 139
$ws.Range('F37').Value = 449
$ws.Range('F38').Value = 72
$ws.Range('F39').Value = 11
$ws.Range('F40').Value = 52
$ws.Range('F41').Value = 2217
$ws.Range('F42').Value = 168
$ws.Range('F43').Value = 115
$ws.Range('F44').Value = 778
$ws.Range('F45').Value = 10
$ws.Range('F46').Value = 133
$ws.Range('F47').Value = 75
$ws.Range('F48').Value = 57
